$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 476.75
$ws.Range("I2").Value = 443.8
$ws.Range("J2").Value = 531.6667
$ws.Range("K2").Value = 443.8
$ws.Range("L2").Value = 531.6667
$ws.Range("M2").Value = -330.8
$ws.Range("N2").Value = -757.6667
$ws.Range("H5").Value = 996.6667
$ws.Range("I5").Value = 184.25
$ws.Range("J5").Value = 2621.5
$ws.Range("K5").Value = 184.25
$ws.Range("L5").Value = 2621.5
$ws.Range("M5").Value = -69.25
$ws.Range("N5").Value = -2851.5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H17").Value = 2077.2173
$ws.Range("J17").Value = 2077.2173
$ws.Range("L17").Value = 6231.651899999999
$ws.Range("N17").Value = -6567.651899999999
$ws.Range("H19").Value = 1247.5927
$ws.Range("I19").Value = 1636.75
$ws.Range("J19").Value = 936.26666
$ws.Range("K19").Value = 1636.75
$ws.Range("L19").Value = 936.26666
$ws.Range("M19").Value = -1461.75
$ws.Range("N19").Value = -1286.26666
$ws.Range("H21").Value = 464999.5
$ws.Range("I21").Value = 464999.5
$ws.Range("K21").Value = 464999.5
$ws.Range("M21").Value = -464531.5
$ws.Range("H23").Value = 464999.5
$ws.Range("I23").Value = 464999.5
$ws.Range("K23").Value = 464999.5
$ws.Range("M23").Value = -464765.5
$ws.Range("H29").Value = 4979
$ws.Range("J29").Value = 4873.75
$ws.Range("L29").Value = 14621.25
$ws.Range("N29").Value = -15183.25
$ws.Range("H38").Value = 10336.091
$ws.Range("J38").Value = 13556.571
$ws.Range("L38").Value = 40669.713
$ws.Range("N38").Value = -41413.713
$ws.Range("H41").Value = 930
$ws.Range("I41").Value = 224
$ws.Range("J41").Value = 2059.6
$ws.Range("K41").Value = 224
$ws.Range("L41").Value = 2059.6
$ws.Range("M41").Value = 216
$ws.Range("N41").Value = -2939.6
$ws.Range("H43").Value = 2997
$ws.Range("I43").Value = 1744.5
$ws.Range("J43").Value = 4249.5
$ws.Range("K43").Value = 1744.5
$ws.Range("L43").Value = 4249.5
$ws.Range("M43").Value = -1675.5
$ws.Range("N43").Value = -4387.5
$ws.Range("H58").Value = 2657
$ws.Range("J58").Value = 2487.5
$ws.Range("L58").Value = 7462.5
$ws.Range("N58").Value = -7762.5
$ws.Range("H64").Value = 80360400
$ws.Range("I64").Value = 62500000
$ws.Range("J64").Value = 83337140
$ws.Range("K64").Value = 62500000
$ws.Range("L64").Value = 83337140
$ws.Range("M64").Value = -62499752
$ws.Range("N64").Value = -83337636
$ws.Range("H67").Value = 80360400
$ws.Range("I67").Value = 62500000
$ws.Range("J67").Value = 83337140
$ws.Range("K67").Value = 62500000
$ws.Range("L67").Value = 83337140
$ws.Range("M67").Value = -62499142
$ws.Range("N67").Value = -83338856
$ws.Range("H87").Value = 77999.57000000001
$ws.Range("J87").Value = 77999.57000000001
$ws.Range("L87").Value = 77999.57000000001
$ws.Range("N87").Value = -80495.57000000001
$ws.Range("H90").Value = 77999.57000000001
$ws.Range("J90").Value = 77999.57000000001
$ws.Range("L90").Value = 233998.71
$ws.Range("N90").Value = -246478.71
$ws.Range("H101").Value = 1520.8
$ws.Range("I101").Value = 1439.1818
$ws.Range("J101").Value = 1745.25
$ws.Range("K101").Value = 4317.5454
$ws.Range("L101").Value = 5235.75
$ws.Range("M101").Value = -2695.5454
$ws.Range("N101").Value = -8479.75
$ws.Range("H107").Value = 1160.9615
$ws.Range("I107").Value = 1195.8695
$ws.Range("J107").Value = 893.3333
$ws.Range("K107").Value = 1195.8695
$ws.Range("L107").Value = 893.3333
$ws.Range("M107").Value = 724.1305
$ws.Range("N107").Value = -4733.3333
$ws.Range("H112").Value = 3274.6223
$ws.Range("J112").Value = 3308.5908
$ws.Range("L112").Value = 9925.7724
$ws.Range("N112").Value = -12141.7724
$ws.Range("H132").Value = 119703.09
$ws.Range("I132").Value = 132201.22
$ws.Range("K132").Value = 396603.66
$ws.Range("M132").Value = -394073.66
$ws.Range("H135").Value = 4178.893
$ws.Range("I135").Value = 1477.4783
$ws.Range("K135").Value = 13297.3047
$ws.Range("M135").Value = -10762.3047
$ws.Range("H136").Value = 107033.71
$ws.Range("J136").Value = 111421.164
$ws.Range("L136").Value = 111421.164
$ws.Range("N136").Value = -121621.164
$ws.Range("I137").Value = 1832.5
$ws.Range("K137").Value = 5497.5
$ws.Range("M137").Value = -2947.5
$ws.Range("H138").Value = 5399.1025
$ws.Range("I138").Value = 4014
$ws.Range("J138").Value = 5535.662
$ws.Range("K138").Value = 12042
$ws.Range("L138").Value = 16606.986
$ws.Range("M138").Value = -6902
$ws.Range("N138").Value = -26886.986

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1152.5625
$ws.Range("I5").Value = 1206.1333
$ws.Range("J5").Value = 349
$ws.Range("K5").Value = 1206.1333
$ws.Range("L5").Value = 349
$ws.Range("M5").Value = -1094.1333
$ws.Range("N5").Value = -573
$ws.Range("H32").Value = 21839.422
$ws.Range("I32").Value = 20432.773
$ws.Range("J32").Value = 28616.908
$ws.Range("K32").Value = 20432.773
$ws.Range("L32").Value = 28616.908
$ws.Range("M32").Value = -20145.773
$ws.Range("N32").Value = -29190.908
$ws.Range("H61").Value = 6522.3945
$ws.Range("I61").Value = 6162.0312
$ws.Range("K61").Value = 6162.0312
$ws.Range("M61").Value = -5950.0312
$ws.Range("H63").Value = 999
$ws.Range("I63").Value = 999
$ws.Range("K63").Value = 999
$ws.Range("M63").Value = -313
$ws.Range("H66").Value = 999
$ws.Range("I66").Value = 999
$ws.Range("K66").Value = 4995
$ws.Range("M66").Value = -1563
$ws.Range("H74").Value = 35716930
$ws.Range("I74").Value = 50002400
$ws.Range("J74").Value = 3250
$ws.Range("K74").Value = 50002400
$ws.Range("L74").Value = 3250
$ws.Range("M74").Value = -50001526
$ws.Range("N74").Value = -4998
$ws.Range("H77").Value = 35716930
$ws.Range("I77").Value = 50002400
$ws.Range("J77").Value = 3250
$ws.Range("K77").Value = 250012000
$ws.Range("L77").Value = 16250
$ws.Range("M77").Value = -250007632
$ws.Range("N77").Value = -24986
$ws.Range("H122").Value = 5044.75
$ws.Range("I122").Value = 4169.9546
$ws.Range("K122").Value = 12509.8638
$ws.Range("M122").Value = -10059.8638
$ws.Range("H132").Value = 11650.934
$ws.Range("I132").Value = 12148.862
$ws.Range("J132").Value = 9952.117
$ws.Range("K132").Value = 36446.586
$ws.Range("L132").Value = 29856.351
$ws.Range("M132").Value = -33916.586
$ws.Range("N132").Value = -34916.351
$ws.Range("H136").Value = 6522.3945
$ws.Range("I136").Value = 6162.0312
$ws.Range("K136").Value = 18486.0936
$ws.Range("M136").Value = -15936.0936
$ws.Range("H138").Value = 77000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 115000
$ws.Range("J140").Value = 115000
$ws.Range("L140").Value = 115000
$ws.Range("N140").Value = -125360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1152.5625
$ws.Range("I4").Value = 1206.1333
$ws.Range("J4").Value = 349
$ws.Range("K4").Value = 1206.1333
$ws.Range("L4").Value = 349
$ws.Range("M4").Value = -1091.1333
$ws.Range("N4").Value = -579
$ws.Range("H22").Value = 403.14285
$ws.Range("I22").Value = 249.36363
$ws.Range("K22").Value = 249.36363
$ws.Range("M22").Value = -76.36363
$ws.Range("H31").Value = 6107.6665
$ws.Range("I31").Value = 4661.5
$ws.Range("J31").Value = 9000
$ws.Range("K31").Value = 4661.5
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = -4409.5
$ws.Range("N31").Value = -9504
$ws.Range("H59").Value = 100097.5
$ws.Range("J59").Value = 100097.5
$ws.Range("L59").Value = 100097.5
$ws.Range("N59").Value = -101791.5
$ws.Range("H68").Value = 50295
$ws.Range("J68").Value = 50295
$ws.Range("L68").Value = 50295
$ws.Range("N68").Value = -51917
$ws.Range("H71").Value = 50295
$ws.Range("J71").Value = 50295
$ws.Range("L71").Value = 150885
$ws.Range("N71").Value = -158997
$ws.Range("H75").Value = 62558.5
$ws.Range("J75").Value = 110117.5
$ws.Range("L75").Value = 110117.5
$ws.Range("N75").Value = -111989.5
$ws.Range("H78").Value = 62558.5
$ws.Range("J78").Value = 110117.5
$ws.Range("L78").Value = 330352.5
$ws.Range("N78").Value = -339712.5
$ws.Range("H86").Value = 4857
$ws.Range("I86").Value = 3764.6667
$ws.Range("J86").Value = 6085.875
$ws.Range("K86").Value = 3764.6667
$ws.Range("L86").Value = 6085.875
$ws.Range("M86").Value = -2641.6667
$ws.Range("N86").Value = -8331.875
$ws.Range("H89").Value = 4857
$ws.Range("I89").Value = 3764.6667
$ws.Range("J89").Value = 6085.875
$ws.Range("K89").Value = 18823.3335
$ws.Range("L89").Value = 30429.375
$ws.Range("M89").Value = -13207.3335
$ws.Range("N89").Value = -41661.375
$ws.Range("H92").Value = 38331.332
$ws.Range("J92").Value = 38331.332
$ws.Range("L92").Value = 38331.332
$ws.Range("N92").Value = -43323.332
$ws.Range("H95").Value = 69597.8
$ws.Range("J95").Value = 69597.8
$ws.Range("L95").Value = 69597.8
$ws.Range("N95").Value = -75089.8
$ws.Range("H97").Value = 31499.834
$ws.Range("J97").Value = 75000
$ws.Range("L97").Value = 75000
$ws.Range("N97").Value = -76982
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 58000
$ws.Range("J102").Value = 84000
$ws.Range("L102").Value = 84000
$ws.Range("N102").Value = -90490
$ws.Range("H103").Value = 57323
$ws.Range("J103").Value = 57323
$ws.Range("L103").Value = 57323
$ws.Range("N103").Value = -59667
$ws.Range("H104").Value = 75000
$ws.Range("J104").Value = 75000
$ws.Range("L104").Value = 75000
$ws.Range("N104").Value = -81988
$ws.Range("H105").Value = 3574.6155
$ws.Range("I105").Value = 3653.7778
$ws.Range("K105").Value = 3653.7778
$ws.Range("M105").Value = -1906.7778
$ws.Range("H106").Value = 33002.5
$ws.Range("J106").Value = 33002.5
$ws.Range("L106").Value = 33002.5
$ws.Range("N106").Value = -35526.5
$ws.Range("H107").Value = 8636.666999999999
$ws.Range("I107").Value = 5955
$ws.Range("K107").Value = 5955
$ws.Range("M107").Value = -4035
$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 80000
$ws.Range("L114").Value = 80000
$ws.Range("N114").Value = -88678
$ws.Range("H116").Value = 83943.39999999999
$ws.Range("J116").Value = 83943.39999999999
$ws.Range("L116").Value = 83943.39999999999
$ws.Range("N116").Value = -93121.39999999999
$ws.Range("H117").Value = 115999
$ws.Range("J117").Value = 115999
$ws.Range("L117").Value = 115999
$ws.Range("N117").Value = -125177
$ws.Range("H119").Value = 30000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 30000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 30000
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -39676
$ws.Range("H120").Value = 51496
$ws.Range("I120").Value = 18000
$ws.Range("J120").Value = 84992
$ws.Range("K120").Value = 18000
$ws.Range("L120").Value = 84992
$ws.Range("M120").Value = -13162
$ws.Range("N120").Value = -94668
$ws.Range("H127").Value = 77994.664
$ws.Range("J127").Value = 77994.664
$ws.Range("L127").Value = 77994.664
$ws.Range("N127").Value = -87914.664
$ws.Range("H129").Value = 49999.8
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
$ws.Range("H134").Value = 3479.2727
$ws.Range("I134").Value = 3665.7778
$ws.Range("K134").Value = 10997.3334
$ws.Range("M134").Value = -8462.3334
$ws.Range("H138").Value = 100099.5
$ws.Range("J138").Value = 100099.5
$ws.Range("L138").Value = 100099.5
$ws.Range("N138").Value = -110379.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 642.7273
$ws.Range("I22").Value = 445.83334
$ws.Range("J22").Value = 879
$ws.Range("K22").Value = 445.83334
$ws.Range("L22").Value = 879
$ws.Range("M22").Value = -95.83334000000002
$ws.Range("N22").Value = -1579
$ws.Range("H31").Value = 21742248
$ws.Range("I31").Value = 30305158
$ws.Range("J31").Value = 5628.769
$ws.Range("K31").Value = 30305158
$ws.Range("L31").Value = 5628.769
$ws.Range("M31").Value = -30304863
$ws.Range("N31").Value = -6218.769
$ws.Range("H34").Value = 21742248
$ws.Range("I34").Value = 30305158
$ws.Range("J34").Value = 5628.769
$ws.Range("K34").Value = 30305158
$ws.Range("L34").Value = 5628.769
$ws.Range("M34").Value = -30304956
$ws.Range("N34").Value = -6032.769
$ws.Range("H97").Value = 17333.166
$ws.Range("J97").Value = 17333.166
$ws.Range("L97").Value = 17333.166
$ws.Range("N97").Value = -19315.166
$ws.Range("H132").Value = 7938072.5
$ws.Range("I132").Value = 9010355
$ws.Range("K132").Value = 27031065
$ws.Range("M132").Value = -27028535

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1391.5
$ws.Range("J2").Value = 2177.4
$ws.Range("L2").Value = 13064.4
$ws.Range("N2").Value = -13290.4
$ws.Range("H4").Value = 4867469.5
$ws.Range("J4").Value = 11835169
$ws.Range("L4").Value = 35505507
$ws.Range("N4").Value = -35505731
$ws.Range("H12").Value = 1039.9231
$ws.Range("J12").Value = 311.25
$ws.Range("L12").Value = 933.75
$ws.Range("N12").Value = -1279.75
$ws.Range("H107").Value = 578.6129
$ws.Range("J107").Value = 685.9
$ws.Range("L107").Value = 2057.7
$ws.Range("N107").Value = -5897.7
$ws.Range("H113").Value = 3608.84
$ws.Range("I113").Value = 1551
$ws.Range("K113").Value = 4653
$ws.Range("M113").Value = -2483
$ws.Range("H131").Value = 7579814
$ws.Range("J131").Value = 10610809
$ws.Range("L131").Value = 31832427
$ws.Range("N131").Value = -31842507
$ws.Range("H140").Value = 4999.8184
$ws.Range("I140").Value = 3894.7368
$ws.Range("J140").Value = 11998.667
$ws.Range("K140").Value = 11684.2104
$ws.Range("L140").Value = 35996.001
$ws.Range("M140").Value = -6504.2104
$ws.Range("N140").Value = -46356.001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 23333500
$ws.Range("I3").Value = 23333500
$ws.Range("K3").Value = 23333500
$ws.Range("M3").Value = -23333384
$ws.Range("H46").Value = 41096.43
$ws.Range("J46").Value = 80000
$ws.Range("L46").Value = 80000
$ws.Range("N46").Value = -80312
$ws.Range("H70").Value = 6313
$ws.Range("I70").Value = 5790
$ws.Range("J70").Value = 7184.6665
$ws.Range("K70").Value = 5790
$ws.Range("L70").Value = 7184.6665
$ws.Range("M70").Value = -5520
$ws.Range("N70").Value = -7724.6665
$ws.Range("H73").Value = 6313
$ws.Range("I73").Value = 5790
$ws.Range("J73").Value = 7184.6665
$ws.Range("K73").Value = 5790
$ws.Range("L73").Value = 7184.6665
$ws.Range("M73").Value = -4854
$ws.Range("N73").Value = -9056.666499999999
$ws.Range("H80").Value = 55664.273
$ws.Range("I80").Value = 117010.336
$ws.Range("J80").Value = 13193.923
$ws.Range("K80").Value = 117010.336
$ws.Range("L80").Value = 13193.923
$ws.Range("M80").Value = -116012.336
$ws.Range("N80").Value = -15189.923
$ws.Range("H83").Value = 55664.273
$ws.Range("I83").Value = 117010.336
$ws.Range("J83").Value = 13193.923
$ws.Range("K83").Value = 585051.6799999999
$ws.Range("L83").Value = 65969.61500000001
$ws.Range("M83").Value = -580059.6799999999
$ws.Range("N83").Value = -75953.61500000001
$ws.Range("H97").Value = 497.5
$ws.Range("I97").Value = 497.5
$ws.Range("K97").Value = 497.5
$ws.Range("M97").Value = -1.5
$ws.Range("H102").Value = 20407320
$ws.Range("I102").Value = 28340528
$ws.Range("K102").Value = 28340528
$ws.Range("M102").Value = -28338906
$ws.Range("H113").Value = 9435709
$ws.Range("I113").Value = 9435709
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 9435709
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -9433539
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 5804.952
$ws.Range("I122").Value = 2007.9565
$ws.Range("J122").Value = 10401.315
$ws.Range("K122").Value = 6023.8695
$ws.Range("L122").Value = 31203.945
$ws.Range("M122").Value = -3573.8695
$ws.Range("N122").Value = -36103.945

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1657666.2
$ws.Range("J20").Value = 7407998
$ws.Range("L20").Value = 7407998
$ws.Range("N20").Value = -7408450
$ws.Range("H22").Value = 1177
$ws.Range("I22").Value = 985.5714
$ws.Range("K22").Value = 985.5714
$ws.Range("M22").Value = -690.5714
$ws.Range("H27").Value = 1177
$ws.Range("I27").Value = 985.5714
$ws.Range("K27").Value = 985.5714
$ws.Range("M27").Value = -878.5714
$ws.Range("H33").Value = 23333.334
$ws.Range("I33").Value = 20000
$ws.Range("J33").Value = 30000
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = -19710
$ws.Range("N33").Value = -30580
$ws.Range("H40").Value = 13895747
$ws.Range("I40").Value = 9622517
$ws.Range("J40").Value = 20839744
$ws.Range("K40").Value = 9622517
$ws.Range("L40").Value = 20839744
$ws.Range("M40").Value = -9622381
$ws.Range("N40").Value = -20840016
$ws.Range("H55").Value = 898.4
$ws.Range("I55").Value = 333.125
$ws.Range("K55").Value = 333.125
$ws.Range("M55").Value = -160.125
$ws.Range("H56").Value = 19577.857
$ws.Range("I56").Value = 10500
$ws.Range("J56").Value = 31681.666
$ws.Range("K56").Value = 10500
$ws.Range("L56").Value = 31681.666
$ws.Range("M56").Value = -9809
$ws.Range("N56").Value = -33063.666
$ws.Range("H61").Value = 3127.1538
$ws.Range("I61").Value = 3189.9473
$ws.Range("J61").Value = 2956.7144
$ws.Range("K61").Value = 3189.9473
$ws.Range("L61").Value = 2956.7144
$ws.Range("M61").Value = -2987.9473
$ws.Range("N61").Value = -3360.7144
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = -21256
$ws.Range("H82").Value = 2267.25
$ws.Range("I82").Value = 957.6
$ws.Range("J82").Value = 3202.7144
$ws.Range("K82").Value = 957.6
$ws.Range("L82").Value = 3202.7144
$ws.Range("M82").Value = -596.6
$ws.Range("N82").Value = -3924.7144
$ws.Range("H85").Value = 2267.25
$ws.Range("I85").Value = 957.6
$ws.Range("J85").Value = 3202.7144
$ws.Range("K85").Value = 957.6
$ws.Range("L85").Value = 3202.7144
$ws.Range("M85").Value = 290.4
$ws.Range("N85").Value = -5698.7144
$ws.Range("H93").Value = 3640.4167
$ws.Range("I93").Value = 3572.5715
$ws.Range("K93").Value = 3572.5715
$ws.Range("M93").Value = -2324.5715
$ws.Range("H100").Value = 2203.0833
$ws.Range("I100").Value = 1762.7142
$ws.Range("K100").Value = 1762.7142
$ws.Range("M100").Value = -1221.7142
$ws.Range("H113").Value = 3127.1538
$ws.Range("I113").Value = 3189.9473
$ws.Range("J113").Value = 2956.7144
$ws.Range("K113").Value = 3189.9473
$ws.Range("L113").Value = 2956.7144
$ws.Range("M113").Value = -1019.9473
$ws.Range("N113").Value = -7296.7144
$ws.Range("H122").Value = 52642532
$ws.Range("I122").Value = 71435520
$ws.Range("J122").Value = 22167.6
$ws.Range("K122").Value = 214306560
$ws.Range("L122").Value = 66502.79999999999
$ws.Range("M122").Value = -214304110
$ws.Range("N122").Value = -71402.79999999999
$ws.Range("H136").Value = 4302.222
$ws.Range("I136").Value = 3890.1428
$ws.Range("J136").Value = 5744.5
$ws.Range("K136").Value = 11670.4284
$ws.Range("L136").Value = 17233.5
$ws.Range("M136").Value = -9120.428400000001
$ws.Range("N136").Value = -22333.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228
$ws.Range("H4").Value = 27983.223
$ws.Range("I4").Value = 50000
$ws.Range("J4").Value = 462.25
$ws.Range("K4").Value = 50000
$ws.Range("L4").Value = 462.25
$ws.Range("M4").Value = -49887
$ws.Range("N4").Value = -688.25
$ws.Range("H37").Value = 36749.668
$ws.Range("J37").Value = 42624.5
$ws.Range("L37").Value = 42624.5
$ws.Range("N37").Value = -43030.5
$ws.Range("H62").Value = 15807.333
$ws.Range("I62").Value = 15198.8
$ws.Range("J62").Value = 16568
$ws.Range("K62").Value = 15198.8
$ws.Range("L62").Value = 16568
$ws.Range("M62").Value = -14574.8
$ws.Range("N62").Value = -17816
$ws.Range("H65").Value = 15807.333
$ws.Range("I65").Value = 15198.8
$ws.Range("J65").Value = 16568
$ws.Range("K65").Value = 75994
$ws.Range("L65").Value = 82840
$ws.Range("M65").Value = -72874
$ws.Range("N65").Value = -89080
$ws.Range("H81").Value = 8945.157999999999
$ws.Range("J81").Value = 13083.4
$ws.Range("L81").Value = 26166.8
$ws.Range("N81").Value = -28288.8
$ws.Range("H84").Value = 8945.157999999999
$ws.Range("J84").Value = 13083.4
$ws.Range("L84").Value = 130834
$ws.Range("N84").Value = -141442
$ws.Range("H88").Value = 33329.332
$ws.Range("J88").Value = 33329.332
$ws.Range("L88").Value = 33329.332
$ws.Range("N88").Value = -34141.332
$ws.Range("H91").Value = 33329.332
$ws.Range("J91").Value = 33329.332
$ws.Range("L91").Value = 33329.332
$ws.Range("N91").Value = -36137.332
$ws.Range("H96").Value = 1169.8334
$ws.Range("I96").Value = 1015.44446
$ws.Range("K96").Value = 1015.44446
$ws.Range("M96").Value = 357.55554
$ws.Range("H100").Value = 1859.95
$ws.Range("I100").Value = 1809.9
$ws.Range("J100").Value = 1910
$ws.Range("K100").Value = 3619.8
$ws.Range("L100").Value = 3820
$ws.Range("M100").Value = -3078.8
$ws.Range("N100").Value = -4902
$ws.Range("H101").Value = 54391.6
$ws.Range("J101").Value = 54391.6
$ws.Range("L101").Value = 54391.6
$ws.Range("N101").Value = -60881.6
$ws.Range("H102").Value = 119999
$ws.Range("J102").Value = 119999
$ws.Range("L102").Value = 119999
$ws.Range("N102").Value = -126489
$ws.Range("H104").Value = 24867.8
$ws.Range("J104").Value = 24867.8
$ws.Range("L104").Value = 24867.8
$ws.Range("N104").Value = -31855.8
$ws.Range("H112").Value = 40461.668
$ws.Range("J112").Value = 40461.668
$ws.Range("L112").Value = 40461.668
$ws.Range("N112").Value = -43415.668
$ws.Range("H113").Value = 598.0909
$ws.Range("I113").Value = 567.9
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 1703.7
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = 466.3000000000002
$ws.Range("N113").Value = -7040
$ws.Range("H114").Value = 119997
$ws.Range("J114").Value = 119997
$ws.Range("L114").Value = 119997
$ws.Range("N114").Value = -128675
$ws.Range("H120").Value = 96000
$ws.Range("J120").Value = 96000
$ws.Range("L120").Value = 96000
$ws.Range("N120").Value = -105676
$ws.Range("H132").Value = 5053123
$ws.Range("I132").Value = 7938202
$ws.Range("J132").Value = 4234.5415
$ws.Range("K132").Value = 23814606
$ws.Range("L132").Value = 12703.6245
$ws.Range("M132").Value = -23812076
$ws.Range("N132").Value = -17763.6245
$ws.Range("H136").Value = 5720.28
$ws.Range("I136").Value = 5125.4585
$ws.Range("K136").Value = 15376.3755
$ws.Range("M136").Value = -12826.3755
